# Update the "想去人数" (F) and "最低票价" (G) figures for the con-list
# entries that changed between scrapes, on both the "展览" and "全部类型"
# worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 2 - F and G both change
    $ws.Range("F2").Value = 7401
    $ws.Range("G2").Value = 62

    # Row 3 - only F changes
    $ws.Range("F3").Value = 67

    # Row 4 - only F changes
    $ws.Range("F4").Value = 209

    # Row 5 - only F changes
    $ws.Range("F5").Value = 209

    # Row 6 - only F changes
    $ws.Range("F6").Value = 1114

    # Row 7 - only F changes
    $ws.Range("F7").Value = 194

    if ($name -eq "展览") {
        # Row 8 - only F changes (this row does not exist in 全部类型's row 8,
        # which instead holds the unrelated "浪漫古典" entry)
        $ws.Range("F8").Value = 13

        # Row 9 - only F changes
        $ws.Range("F9").Value = 107

        # Row 10 - only F changes
        $ws.Range("F10").Value = 28
    } else {
        # On 全部类型 the rows are shifted down by one starting at row 9
        # because of the extra "浪漫古典" row inserted at row 8.

        # Row 9 - only F changes
        $ws.Range("F9").Value = 13

        # Row 10 - only F changes
        $ws.Range("F10").Value = 107

        # Row 11 - only F changes
        $ws.Range("F11").Value = 28
    }
}
